$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 1 (the "ID" row), shifting everything up
$ws.Rows("1").Delete()

# Delete column A (now empty since the old "ID" cell was removed), shifting B->A, C->B
$ws.Columns("A").Delete()

Write-Host "A1=" $ws.Range("A1").Value2
Write-Host "B1=" $ws.Range("B1").Value2
Write-Host "A2=" $ws.Range("A2").Value2
Write-Host "B2=" $ws.Range("B2").Value2
Write-Host "A3=" $ws.Range("A3").Value2
Write-Host "B3=" $ws.Range("B3").Value2
Write-Host "dim=" $ws.UsedRange.Address()
